$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''65.721.32'
$ws.Range('E2').Value = '  +1.14%  '
$ws.Range('D3').Value = '''3.305.75'
$ws.Range('E3').Value = '  +1.50%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '''186.96'
$ws.Range('E5').Value = '  +5.77%  '
$ws.Range('D6').Value = '''552.90'
$ws.Range('E6').Value = '  +0.58%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D8').Value = '''3.300.12'
$ws.Range('E8').Value = '  +1.63%  '
$ws.Range('B9').Value = 'XRP'
$ws.Range('C9').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D9').Value = '''0.578'
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('D10').Value = '''0.178'
$ws.Range('E10').Value = '  -2.64%  '
$ws.Range('D11').Value = '''0.579'
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('D12').Value = '''46.56'
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('D13').Value = '''0.0000265'
$ws.Range('E13').Value = '  +2.36%  '
$ws.Range('D14').Value = '''8.58'
$ws.Range('E14').Value = '  +1.83%  '
$ws.Range('D15').Value = '''3.836.30'
$ws.Range('E15').Value = '  +1.47%  '
$ws.Range('D16').Value = '''595.97'
$ws.Range('E16').Value = '  -0.75%  '
$ws.Range('D17').Value = '''65.792.81'
$ws.Range('E17').Value = '  +1.47%  '
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('D19').Value = '''17.84'
$ws.Range('E19').Value = '  +0.61%  '
$ws.Range('D20').Value = '''3.302.87'
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('D21').Value = '''10.97'
$ws.Range('E21').Value = '  -2.57%  '
$ws.Range('D22').Value = '''0.894'
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('D23').Value = '''18.40'
$ws.Range('E23').Value = '  +6.11%  '
$ws.Range('D24').Value = '''5.06'
$ws.Range('E24').Value = '  +3.09%  '
$ws.Range('D25').Value = '''100.26'
$ws.Range('E25').Value = '  -1.07%  '
$ws.Range('D26').Value = '''3.92'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('E27').Value = '  +1.52%  '
$ws.Range('D28').Value = '''2.72'
$ws.Range('E28').Value = '  +3.21%  '
$ws.Range('D29').Value = '''9.42'
$ws.Range('E29').Value = '  +2.63%  '
$ws.Range('D30').Value = '''8.64'
$ws.Range('E30').Value = '  +1.19%  '
$ws.Range('D31').Value = '''30.27'
$ws.Range('E31').Value = '  +0.72%  '
$ws.Range('D32').Value = '''6.66'
$ws.Range('E32').Value = '  +8.39%  '
$ws.Range('D33').Value = '''3.80'
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('D34').Value = '''568.73'
$ws.Range('E34').Value = '  +7.94%  '
$ws.Range('D35').Value = '''10.96'
$ws.Range('E35').Value = '  +0.70%  '
$ws.Range('D36').Value = '''0.103'
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').Value = '''0.999'
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '''3.702.90'
$ws.Range('E38').Value = '  -0.78%  '
$ws.Range('D39').Value = '''56.66'
$ws.Range('E39').Value = '  +1.79%  '
$ws.Range('D40').Value = '''3.45'
$ws.Range('E40').Value = '  +9.10%  '
$ws.Range('D41').Value = '''33.50'
$ws.Range('E41').Value = '  +6.63%  '
$ws.Range('D42').Value = '''3.21'
$ws.Range('E42').Value = '  -5.59%  '
$ws.Range('E43').Value = '  +2.63%  '
$ws.Range('B44').Value = 'PEPE'
$ws.Range('C44').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D44').Value = '''0.0₃0693'
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').Value = '''3.42'
$ws.Range('E45').Value = '  +8.10%  '
$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D46').Value = '''2.63'
$ws.Range('E46').Value = '  +0.22%  '
$ws.Range('D47').Value = '''0.335'
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('D48').Value = '''0.0415'
$ws.Range('E48').Value = '  +2.81%  '
$ws.Range('D49').Value = '''0.128'
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('B50').Value = 'FirstDigitalUSD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D50').Value = '''1.00'
$ws.Range('E50').Value = '  +0.30%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').Value = '''2.54'
$ws.Range('E51').Value = '  -0.41%  '
